# Generate Report for Handoff
# Updates the "Priority" column (E) for newly-handed-off rows to "ht",
# and refreshes the handoff timestamps on the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 13, 14)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

foreach ($r in $rows) {
    # Priority column (E) on the per-language sheets: empty -> "ht"
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"

    # Latest HO Xliff Generate Date on Overview (column G)
    $wsOverview.Range("G$r").Value = "2016-08-30 08:23:53"

    # Latest Handoff Datetime on zh-cn (column H)
    $wsZhCn.Range("H$r").Value = "2016-08-30 08:23:48"

    # Latest Handoff Datetime on de-de (column H)
    $wsDeDe.Range("H$r").Value = "2016-08-30 08:23:53"
}
